$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.578493666666667
$ws.Range("H2").Value = 16.735481
$ws.Range("I2").Value = 0.1036332930693284
$ws.Range("J2").Value = 0.1036332930693284
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08343099999999999
$ws.Range("N2").Value = 0.250293
$ws.Range("O2").Value = 0.0134664339880133
$ws.Range("P2").Value = 0.0134664339880133
$ws.Range("Q2").Value = 0.4654193051036666
$ws.Range("R2").Value = 4.188773745933
$ws.Range("S2").Value = 0.001395570900078547
$ws.Range("T2").Value = 0.001395570900078547
$ws.Range("G3").Value = 5.578493666666667
$ws.Range("H3").Value = 16.735481
$ws.Range("I3").Value = 0.1036332930693284
$ws.Range("J3").Value = 0.1036332930693284
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.699187666666667
$ws.Range("N3").Value = 14.097563
$ws.Range("O3").Value = 0.758486659760196
$ws.Range("P3").Value = 0.758486659760196
$ws.Range("Q3").Value = 26.21438863697811
$ws.Range("R3").Value = 235.929497732803
$ws.Range("S3").Value = 0.07860447030010437
$ws.Range("T3").Value = 0.07860447030010438
$ws.Range("G4").Value = 5.578493666666667
$ws.Range("H4").Value = 16.735481
$ws.Range("I4").Value = 0.1036332930693284
$ws.Range("J4").Value = 0.1036332930693284
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1492686666666667
$ws.Range("N4").Value = 0.447806
$ws.Range("O4").Value = 0.02409316256721636
$ws.Range("P4").Value = 0.02409316256721636
$ws.Range("Q4").Value = 0.8326943116317779
$ws.Range("R4").Value = 7.494248804686
$ws.Range("S4").Value = 0.002496853777295306
$ws.Range("T4").Value = 0.002496853777295306
$ws.Range("G5").Value = 5.578493666666667
$ws.Range("H5").Value = 16.735481
$ws.Range("I5").Value = 0.1036332930693284
$ws.Range("J5").Value = 0.1036332930693284
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.263591
$ws.Range("N5").Value = 3.790773
$ws.Range("O5").Value = 0.2039537436845743
$ws.Range("P5").Value = 0.2039537436845743
$ws.Range("Q5").Value = 7.048934390757001
$ws.Range("R5").Value = 63.44040951681301
$ws.Range("S5").Value = 0.02113639809185017
$ws.Range("T5").Value = 0.02113639809185017
$ws.Range("G6").Value = 3.704475666666667
$ws.Range("H6").Value = 11.113427
$ws.Range("I6").Value = 0.06881911773528272
$ws.Range("J6").Value = 0.06881911773528274
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.08343099999999999
$ws.Range("N6").Value = 0.250293
$ws.Range("O6").Value = 0.0134664339880133
$ws.Range("P6").Value = 0.0134664339880133
$ws.Range("Q6").Value = 0.3090681093456666
$ws.Range("R6").Value = 2.781612984111
$ws.Range("S6").Value = 0.0009267481060955
$ws.Range("T6").Value = 0.0009267481060955002
$ws.Range("G7").Value = 3.704475666666667
$ws.Range("H7").Value = 11.113427
$ws.Range("I7").Value = 0.06881911773528272
$ws.Range("J7").Value = 0.06881911773528274
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.699187666666667
$ws.Range("N7").Value = 14.097563
$ws.Range("O7").Value = 0.758486659760196
$ws.Range("P7").Value = 0.758486659760196
$ws.Range("Q7").Value = 17.40802636426678
$ws.Range("R7").Value = 156.672237278401
$ws.Range("S7").Value = 0.05219838273867826
$ws.Range("T7").Value = 0.05219838273867827
$ws.Range("G8").Value = 3.704475666666667
$ws.Range("H8").Value = 11.113427
$ws.Range("I8").Value = 0.06881911773528272
$ws.Range("J8").Value = 0.06881911773528274
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1492686666666667
$ws.Range("N8").Value = 0.447806
$ws.Range("O8").Value = 0.02409316256721636
$ws.Range("P8").Value = 0.02409316256721636
$ws.Range("Q8").Value = 0.5529621434624445
$ws.Range("R8").Value = 4.976659291162
$ws.Range("S8").Value = 0.001658070191328569
$ws.Range("T8").Value = 0.00165807019132857
$ws.Range("G9").Value = 3.704475666666667
$ws.Range("H9").Value = 11.113427
$ws.Range("I9").Value = 0.06881911773528272
$ws.Range("J9").Value = 0.06881911773528274
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.263591
$ws.Range("N9").Value = 3.790773
$ws.Range("O9").Value = 0.2039537436845743
$ws.Range("P9").Value = 0.2039537436845743
$ws.Range("Q9").Value = 4.680942112119
$ws.Range("R9").Value = 42.128479009071
$ws.Range("S9").Value = 0.01403591669918039
$ws.Range("T9").Value = 0.01403591669918039
$ws.Range("G10").Value = 24.77295966666667
$ws.Range("H10").Value = 74.31887900000001
$ws.Range("I10").Value = 0.4602144490493554
$ws.Range("J10").Value = 0.4602144490493556
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.08343099999999999
$ws.Range("N10").Value = 0.250293
$ws.Range("O10").Value = 0.0134664339880133
$ws.Range("P10").Value = 0.0134664339880133
$ws.Range("Q10").Value = 2.066832797949667
$ws.Range("R10").Value = 18.601495181547
$ws.Range("S10").Value = 0.006197447498453054
$ws.Range("T10").Value = 0.006197447498453058
$ws.Range("G11").Value = 24.77295966666667
$ws.Range("H11").Value = 74.31887900000001
$ws.Range("I11").Value = 0.4602144490493554
$ws.Range("J11").Value = 0.4602144490493556
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.699187666666667
$ws.Range("N11").Value = 14.097563
$ws.Range("O11").Value = 0.758486659760196
$ws.Range("P11").Value = 0.758486659760196
$ws.Range("Q11").Value = 116.4127865324308
$ws.Range("R11").Value = 1047.715078791877
$ws.Range("S11").Value = 0.3490665202328245
$ws.Range("T11").Value = 0.3490665202328246
$ws.Range("G12").Value = 24.77295966666667
$ws.Range("H12").Value = 74.31887900000001
$ws.Range("I12").Value = 0.4602144490493554
$ws.Range("J12").Value = 0.4602144490493556
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1492686666666667
$ws.Range("N12").Value = 0.447806
$ws.Range("O12").Value = 0.02409316256721636
$ws.Range("P12").Value = 0.02409316256721636
$ws.Range("Q12").Value = 3.697826658830445
$ws.Range("R12").Value = 33.28043992947401
$ws.Range("S12").Value = 0.01108802153672803
$ws.Range("T12").Value = 0.01108802153672804
$ws.Range("G13").Value = 24.77295966666667
$ws.Range("H13").Value = 74.31887900000001
$ws.Range("I13").Value = 0.4602144490493554
$ws.Range("J13").Value = 0.4602144490493556
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.263591
$ws.Range("N13").Value = 3.790773
$ws.Range("O13").Value = 0.2039537436845743
$ws.Range("P13").Value = 0.2039537436845743
$ws.Range("Q13").Value = 31.30288887816301
$ws.Range("R13").Value = 281.725999903467
$ws.Range("S13").Value = 0.0938624597813498
$ws.Range("T13").Value = 0.09386245978134983
$ws.Range("G14").Value = 0.4291063333333334
$ws.Range("H14").Value = 1.287319
$ws.Range("I14").Value = 0.007971632676749163
$ws.Range("J14").Value = 0.007971632676749165
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.08343099999999999
$ws.Range("N14").Value = 0.250293
$ws.Range("O14").Value = 0.0134664339880133
$ws.Range("P14").Value = 0.0134664339880133
$ws.Range("Q14").Value = 0.03580077049633333
$ws.Range("R14").Value = 0.322206934467
$ws.Range("S14").Value = 0.0001073494652181323
$ws.Range("T14").Value = 0.0001073494652181324
$ws.Range("G15").Value = 0.4291063333333334
$ws.Range("H15").Value = 1.287319
$ws.Range("I15").Value = 0.007971632676749163
$ws.Range("J15").Value = 0.007971632676749165
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.699187666666667
$ws.Range("N15").Value = 14.097563
$ws.Range("O15").Value = 0.758486659760196
$ws.Range("P15").Value = 0.758486659760196
$ws.Range("Q15").Value = 2.016451189288556
$ws.Range("R15").Value = 18.148060703597
$ws.Range("S15").Value = 0.006046377041822703
$ws.Range("T15").Value = 0.006046377041822705
$ws.Range("G16").Value = 0.4291063333333334
$ws.Range("H16").Value = 1.287319
$ws.Range("I16").Value = 0.007971632676749163
$ws.Range("J16").Value = 0.007971632676749165
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1492686666666667
$ws.Range("N16").Value = 0.447806
$ws.Range("O16").Value = 0.02409316256721636
$ws.Range("P16").Value = 0.02409316256721636
$ws.Range("Q16").Value = 0.0640521302348889
$ws.Range("R16").Value = 0.5764691721140001
$ws.Range("S16").Value = 0.0001920618420070517
$ws.Range("T16").Value = 0.0001920618420070517
$ws.Range("G17").Value = 0.4291063333333334
$ws.Range("H17").Value = 1.287319
$ws.Range("I17").Value = 0.007971632676749163
$ws.Range("J17").Value = 0.007971632676749165
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.263591
$ws.Range("N17").Value = 3.790773
$ws.Range("O17").Value = 0.2039537436845743
$ws.Range("P17").Value = 0.2039537436845743
$ws.Range("Q17").Value = 0.5422149008430001
$ws.Range("R17").Value = 4.879934107587
$ws.Range("S17").Value = 0.001625844327701275
$ws.Range("T17").Value = 0.001625844327701276
$ws.Range("G18").Value = 19.34413
$ws.Range("H18").Value = 58.03239
$ws.Range("I18").Value = 0.3593615074692841
$ws.Range("J18").Value = 0.3593615074692842
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.08343099999999999
$ws.Range("N18").Value = 0.250293
$ws.Range("O18").Value = 0.0134664339880133
$ws.Range("P18").Value = 0.0134664339880133
$ws.Range("Q18").Value = 1.61390011003
$ws.Range("R18").Value = 14.52510099027
$ws.Range("S18").Value = 0.004839318018168062
$ws.Range("T18").Value = 0.004839318018168064
$ws.Range("G19").Value = 19.34413
$ws.Range("H19").Value = 58.03239
$ws.Range("I19").Value = 0.3593615074692841
$ws.Range("J19").Value = 0.3593615074692842
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 4.699187666666667
$ws.Range("N19").Value = 14.097563
$ws.Range("O19").Value = 0.758486659760196
$ws.Range("P19").Value = 0.758486659760196
$ws.Range("Q19").Value = 90.90169711839667
$ws.Range("R19").Value = 818.11527406557
$ws.Range("S19").Value = 0.272570909446766
$ws.Range("T19").Value = 0.2725709094467661
$ws.Range("G20").Value = 19.34413
$ws.Range("H20").Value = 58.03239
$ws.Range("I20").Value = 0.3593615074692841
$ws.Range("J20").Value = 0.3593615074692842
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.1492686666666667
$ws.Range("N20").Value = 0.447806
$ws.Range("O20").Value = 0.02409316256721636
$ws.Range("P20").Value = 0.02409316256721636
$ws.Range("Q20").Value = 2.887472492926667
$ws.Range("R20").Value = 25.98725243634
$ws.Range("S20").Value = 0.0086581552198574
$ws.Range("T20").Value = 0.0086581552198574
$ws.Range("G21").Value = 19.34413
$ws.Range("H21").Value = 58.03239
$ws.Range("I21").Value = 0.3593615074692841
$ws.Range("J21").Value = 0.3593615074692842
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 1.263591
$ws.Range("N21").Value = 3.790773
$ws.Range("O21").Value = 0.2039537436845743
$ws.Range("P21").Value = 0.2039537436845743
$ws.Range("Q21").Value = 24.44306857083
$ws.Range("R21").Value = 219.98761713747
$ws.Range("S21").Value = 0.0732931247844926
$ws.Range("T21").Value = 0.0732931247844926
